# Update the quote date (A1) by one day, and refresh the unit prices
# for the three hinge ("Bisagra") items (step 1 and 2 price fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45309

$ws.Range("D26").Value = 338.256
$ws.Range("D27").Value = 389.76
$ws.Range("D28").Value = 478.5
